# Auto-generated Word COM-interop script implementing the target diff.
$d = $word.ActiveDocument

# --- Hunk 1: "Objectives" bullet -> AI-Powered Intent Recognition --------------
$null = $d.Content.Find.Execute(
    "To implement intelligent tool suggestion using natural language processing techniques (TF-IDF and cosine similarity) based on user descriptions.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "AI-Powered Intent Recognition: Implement state-of-the-art semantic understanding using sentence transformers to accurately match user descriptions with appropriate image effects.​", 2)

# --- Hunk 2 & 3: End-user bullets get bold names + 4 new persona bullets -------
# Locate the "General Users" and "Design Enthusiasts" bullet paragraphs by their
# (still-unmodified at this point) text, then replace that whole range with the
# fully-specified OOXML for the updated/new paragraphs in one shot.
$generalUsersPara = $null
$designEnthusiastsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text.StartsWith("General Users:")) {
        $generalUsersPara = $para
    } elseif ($text.StartsWith("Design Enthusiasts:")) {
        $designEnthusiastsPara = $para
    }
}
$spanStart = $generalUsersPara.Range.Start
$spanEnd = $designEnthusiastsPara.Range.End
$span = $d.Range($spanStart, $spanEnd)
$span.InsertXML("<w:p><w:pPr><w:pStyle w:val=`"ListBullet`"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>General Users</w:t></w:r><w:r><w:t xml:space=`"preserve`"> -</w:t></w:r><w:r><w:t xml:space=`"preserve`"> Individuals interested in experimenting with image effects for creative or fun purposes, such as hobbyists or students.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListBullet`"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Design Enthusiasts</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`"> and</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t>Digital Artists</w:t></w:r><w:r><w:t xml:space=`"preserve`"> - </w:t></w:r><w:r><w:t>Users who want quick transformations like sepia, sketches, or pixelation without complex software like Photoshop.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:rPr><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t>Content Creators</w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t xml:space=`"preserve`"> - </w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t>Social media influencers, bloggers, and digital artists seeking quick aesthetic transformations for their visual content with minimal friction.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:rPr><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t>Students and Educators</w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t xml:space=`"preserve`"> - </w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t>Academic users exploring image processing concepts through hands-on experimentation with various filters and effects.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:rPr><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t>Hobbyists</w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t xml:space=`"preserve`"> - </w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-IN`" w:eastAsia=`"en-IN`"/></w:rPr><w:t>Photography enthusiasts experimenting with artistic transformations and retro aesthetics for personal projects.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"0`"/></w:numPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr></w:p>")

# --- Hunk 4: drop the stale lastRenderedPageBreak on the "Aesthetic Blur" bullet
# A benign find/replace (old text -> identical text) forces Word to re-serialize
# the run, which drops the cached lastRenderedPageBreak marker.
$null = $d.Content.Find.Execute(
    "  Aesthetic Blur: Applies Gaussian blur for a soft, dreamy effect.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "  Aesthetic Blur: Applies Gaussian blur for a soft, dreamy effect.", 2)

# --- Hunk 5: register the new "my-2" paragraph style ---------------------------
$my2 = $d.Styles.Add("my-2", 1)
$my2.BaseStyle = $d.Styles.Item("Normal")
$my2.Font.Name = "Times New Roman"
$my2.Font.NameFarEast = "Times New Roman"
$my2.Font.NameBi = "Times New Roman"
$my2.Font.Size = 12
$my2.Font.SizeBi = 12
$my2.Font.LanguageID = "en-IN"
$my2.Font.LanguageIDFarEast = "en-IN"
$my2.ParagraphFormat.SpaceBefore = 5
$my2.ParagraphFormat.SpaceBeforeAuto = $true
$my2.ParagraphFormat.SpaceAfter = 5
$my2.ParagraphFormat.SpaceAfterAuto = $true
$my2.ParagraphFormat.LineSpacingRule = 0

Write-Output "edit complete"
